$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -0.02283115780924631
$ws.Range("E4").Value = 0.03252158606886344
$ws.Range("J4").Value = -0.01001126144342611
$ws.Range("C5").Value = 0.007191296927651877
$ws.Range("E5").Value = 0.0294076926963077
$ws.Range("J5").Value = 0.009471818331361903
$ws.Range("C6").Value = 0.1519870517114821
$ws.Range("E6").Value = 0.9382294163451764
$ws.Range("J6").Value = -0.002172398423935881
$ws.Range("C7").Value = 0.9484030940321236
$ws.Range("E7").Value = -0.01072744132509765
$ws.Range("J7").Value = 0.007057077690483055
$ws.Range("C8").Value = 0.0225188151727526
$ws.Range("E8").Value = 0.01356077353443094
$ws.Range("J8").Value = 0.01109177598846261
$ws.Range("C9").Value = -0.004604815384192615
$ws.Range("E9").Value = -0.007437114921484595
$ws.Range("J9").Value = -0.01220219018234914
$ws.Range("C10").Value = 0.04784935650597425
$ws.Range("E10").Value = 0.02169796397191855
$ws.Range("J10").Value = 0.006567211724689074
$ws.Range("C11").Value = 0.1444172317606893
$ws.Range("E11").Value = 0.003330341797213672
$ws.Range("J11").Value = -0.005271696816414461
$ws.Range("C12").Value = -0.1937698121987924
$ws.Range("E12").Value = 0.01269399737175989
$ws.Range("J12").Value = -0.003221754185824128
$ws.Range("C13").Value = 0.02240781219231248
$ws.Range("E13").Value = -0.001322413396896536
$ws.Range("J13").Value = 0.01383563057094249
$ws.Range("C14").Value = 0.00005901744236069768
$ws.Range("E14").Value = -0.005465977850639113
$ws.Range("J14").Value = -0.0006871099755667106
$ws.Range("C15").Value = -0.01752631202905248
$ws.Range("E15").Value = -0.007758804982352198
$ws.Range("J15").Value = 0.02159517101128804
$ws.Range("C16").Value = 0.001620135616805425
$ws.Range("E16").Value = 0.007825959385038373
$ws.Range("J16").Value = -0.01956840203256883
$ws.Range("C17").Value = 0.00136514549460582
$ws.Range("E17").Value = 0.000007960320318412811
$ws.Range("J17").Value = -0.004315955290631166
$ws.Range("C18").Value = 0.01526795734671829
$ws.Range("E18").Value = 0.002708669676346787
$ws.Range("J18").Value = 0.007453877028566195
$ws.Range("C19").Value = 0.005928089613123583
$ws.Range("E19").Value = 0.0006603959304158371
$ws.Range("J19").Value = -0.009818421305930008
$ws.Range("C20").Value = -0.0181484401499376
$ws.Range("E20").Value = -0.006156371478254857
$ws.Range("J20").Value = 0.003543377342157545
$ws.Range("C21").Value = -0.022813675152547
$ws.Range("E21").Value = 0.0113808977032359
$ws.Range("J21").Value = 0.008196569064025674
$ws.Range("C22").Value = -0.02214966990998679
$ws.Range("E22").Value = -0.01159324654372986
$ws.Range("J22").Value = 0.006366690229209365
$ws.Range("C23").Value = -0.004630091513203659
$ws.Range("E23").Value = -0.01230116593204664
$ws.Range("J23").Value = -0.01849373487128882
$ws.Range("C24").Value = -0.009812086472483458
$ws.Range("E24").Value = 0.007106797916271916
$ws.Range("J24").Value = 0.004913271513563789
$ws.Range("C25").Value = 0.004657418106296723
$ws.Range("E25").Value = 0.01962831332913253
$ws.Range("J25").Value = 0.005847604250277086
$ws.Range("C26").Value = 0.02198199620727985
$ws.Range("E26").Value = 0.005956423246256929
$ws.Range("J26").Value = -0.0129405840538354
$ws.Range("C27").Value = 0.02228941845957674
$ws.Range("E27").Value = -0.02877250838290033
$ws.Range("J27").Value = 0.001217981061547961
$ws.Range("C28").Value = -0.01376402301456092
$ws.Range("E28").Value = 0.01687623609904944
$ws.Range("J28").Value = -0.01671249160560131
$ws.Range("C29").Value = 0.01106015948240638
$ws.Range("E29").Value = -0.005511798460471938
$ws.Range("J29").Value = 0.01089234380240115
$ws.Range("C30").Value = 0.003827297721091908
$ws.Range("E30").Value = 0.005009897960395917
$ws.Range("J30").Value = -0.021808116186103
$ws.Range("C31").Value = -0.00755999828639993
$ws.Range("E31").Value = 0.01726961963478478
$ws.Range("J31").Value = 0.01263127517042416
$ws.Range("C32").Value = 0.008916752324670093
$ws.Range("E32").Value = 0.003597134351885373
$ws.Range("J32").Value = 0.0007849565781093693
$ws.Range("C33").Value = 0.006172791318911653
$ws.Range("E33").Value = -0.001873426922937076
$ws.Range("J33").Value = 0.002380971712721403
